# Add a new worksheet "AddLeave" after the existing "Sheet1" and populate it
# with an employee leave-entitlement table (mirrors reading multiple
# sheets/rows from the HR automation workbook).

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws.Name = "AddLeave"

# Header row
$ws.Range("A1").Value = "Employee"
$ws.Range("B1").Value = "Leave Type"
$ws.Range("C1").Value = "Leave Period"
$ws.Range("D1").Value = "Entitlement"
$ws.Range("A1:D1").Font.Bold = $true

# Make sure the data cells are stored as text (numbers/dates keep their
# literal textual form, e.g. "2020-01-01$$2020-12-31").
$ws.Range("A2:D3").NumberFormat = "@"

# Row 2 - filled in field order: employee, leave period, leave type, entitlement
$ws.Range("A2").Value = "Russel"
$ws.Range("C2").Value = "2020-01-01`$`$2020-12-31"
$ws.Range("B2").Value = "4"
$ws.Range("D2").Value = "20"

# Row 3
$ws.Range("A3").Value = "Russel"
$ws.Range("C3").Value = "2020-01-01`$`$2020-12-31"
$ws.Range("B3").Value = "3"
$ws.Range("D3").Value = "30"

# Widen the "Leave Period" column so the date-range text is fully visible.
$ws.Columns("C").ColumnWidth = 21.5

# Match the author's last on-screen state of the sheet.
$ws.Application.ActiveWindow.Zoom = 223
$ws.Range("C6").Select() | Out-Null
